$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33 currently holds the "Charles University in Prague" / "A course on
# Spatial Data Science" / "https://martinfleischmann.net/sds/" / "GeoData"
# entry. The title becomes more specific, and a new row is inserted right
# below it for the related "micro" course with the same author/link pattern.

# 1) Rename the existing course title on row 33.
$ws.Range("B33").Value = "A course on Spatial Data Science for Social Geography"

# 2) Insert a new row right after row 33 for the related course.
$ws.Rows.Item(34).Insert()

$ws.Range("A34").Value = "Charles University in Prague"
$ws.Range("B34").Value = "A course on Spatial Data Science"
$ws.Range("C34").Value = "https://martinfleischmann.net/sds/micro/"
$ws.Range("D34").Value = "GeoData"

# Match the wrapped-text row height used by the other multi-line rows.
$ws.Rows.Item(34).RowHeight = 34

# Update the visible selection to reflect the newly edited/inserted rows.
$ws.Range("A33:D34").Select()
